$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price/Volume columns for the rows being
# updated so numeric-looking strings (e.g. "593.32") are stored as text,
# matching the inlineStr cells in the original workbook, then restore the
# default "Normal" style so no stray style index is left on the cells.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = "70.994.65"
$ws.Range("E2").Value = "  +5.69%  "
$ws.Range("D3").Value = "3.645.77"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "593.32"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").Value = "195.07"
$ws.Range("E6").Value = "  +3.70%  "
$ws.Range("D7").Value = "0.644"
$ws.Range("E7").Value = "  +2.22%  "
$ws.Range("D8").Value = "3.640.80"
$ws.Range("E8").Value = "  +5.61%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "0.185"
$ws.Range("E10").Value = "  +7.22%  "
$ws.Range("D11").Value = "0.677"
$ws.Range("E11").Value = "  +5.01%  "
$ws.Range("D12").Value = "58.00"
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "0.0000311"
$ws.Range("E13").Value = "  +12.24%  "
$ws.Range("D14").Value = "9.94"
$ws.Range("E14").Value = "  +4.76%  "
$ws.Range("D15").Value = "4.230.53"
$ws.Range("E15").Value = "  +5.97%  "
$ws.Range("D16").Value = "20.61"
$ws.Range("E16").Value = "  +8.70%  "
$ws.Range("D17").Value = "3.644.67"
$ws.Range("E17").Value = "  +5.62%  "
$ws.Range("D18").Value = "70.938.01"
$ws.Range("D19").Value = "12.78"
$ws.Range("E19").Value = "  +5.90%  "
$ws.Range("E20").Value = "  +3.04%  "
$ws.Range("E21").Value = "  +4.09%  "
$ws.Range("D22").Value = "492.07"
$ws.Range("E22").Value = "  +2.32%  "
$ws.Range("D23").Value = "18.92"
$ws.Range("E23").Value = "  +8.95%  "
$ws.Range("D24").Value = "5.23"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "4.54"
$ws.Range("E25").Value = "  +4.40%  "
$ws.Range("D26").Value = "91.16"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("D27").Value = "3.17"
$ws.Range("E27").Value = "  +6.80%  "
$ws.Range("D28").Value = "11.52"
$ws.Range("E28").Value = "  +4.78%  "
$ws.Range("D29").Value = "9.60"
$ws.Range("E29").Value = "  +6.45%  "
$ws.Range("D30").Value = "7.95"
$ws.Range("E30").Value = "  +7.67%  "
$ws.Range("D31").Value = "32.83"
$ws.Range("E31").Value = "  +4.96%  "
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").Value = "  +8.94%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "68.13"
$ws.Range("E33").Value = "  +4.84%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "12.30"
$ws.Range("E34").Value = "  +4.11%  "
$ws.Range("D35").Value = "616.01"
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("D36").Value = "40.37"
$ws.Range("E36").Value = "  +9.02%  "
$ws.Range("D37").Value = "0.0₃0841"
$ws.Range("E37").Value = "  +8.77%  "
$ws.Range("D38").Value = "0.412"
$ws.Range("E38").Value = "  +6.37%  "
$ws.Range("D39").Value = "0.148"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("D41").Value = "3.58"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("D42").Value = "3.323.56"
$ws.Range("E42").Value = "  +3.80%  "
$ws.Range("D43").Value = "2.94"
$ws.Range("E43").Value = "  +14.26%  "
$ws.Range("D44").Value = "3.23"
$ws.Range("E44").Value = "  +19.87%  "
$ws.Range("E45").Value = "  +9.79%  "
$ws.Range("D46").Value = "0.0461"
$ws.Range("E46").Value = "  +7.00%  "
$ws.Range("D47").Value = "9.70"
$ws.Range("E47").Value = "  +11.96%  "
$ws.Range("D48").Value = "3.34"
$ws.Range("E48").Value = "  +3.59%  "
$ws.Range("D49").Value = "0.140"
$ws.Range("E49").Value = "  +3.16%  "
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "3.23"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  -0.05%  "

$rng.Style = "Normal"
